# Auto-generated update of market price columns (H-N) across all job sheets
# Mirrors a scheduled Moogle market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 753.6667
$ws.Range("I18").Value = 753.6667
$ws.Range("K18").Value = 753.6667
$ws.Range("M18").Value = -469.6667
$ws.Range("H40").Value = 9882.823
$ws.Range("I40").Value = 9409.143
$ws.Range("J40").Value = 10214.4
$ws.Range("K40").Value = 9409.143
$ws.Range("L40").Value = 10214.4
$ws.Range("M40").Value = -9234.143
$ws.Range("N40").Value = -10564.4
$ws.Range("H70").Value = 5239.5454
$ws.Range("J70").Value = 6079.375
$ws.Range("L70").Value = 18238.125
$ws.Range("N70").Value = -18778.125
$ws.Range("H73").Value = 5239.5454
$ws.Range("J73").Value = 6079.375
$ws.Range("L73").Value = 18238.125
$ws.Range("N73").Value = -20110.125
$ws.Range("H76").Value = 6963.4287
$ws.Range("I76").Value = 5998.1665
$ws.Range("K76").Value = 5998.1665
$ws.Range("M76").Value = -5683.1665
$ws.Range("H79").Value = 6963.4287
$ws.Range("I79").Value = 5998.1665
$ws.Range("K79").Value = 5998.1665
$ws.Range("M79").Value = -4906.1665
$ws.Range("H96").Value = 336.4
$ws.Range("I96").Value = 336.4
$ws.Range("K96").Value = 1009.2
$ws.Range("M96").Value = 363.8000000000001
$ws.Range("H106").Value = 62861180
$ws.Range("I106").Value = 146670640
$ws.Range("J106").Value = 4089.25
$ws.Range("K106").Value = 146670640
$ws.Range("L106").Value = 4089.25
$ws.Range("M106").Value = -146670009
$ws.Range("N106").Value = -5351.25
$ws.Range("H113").Value = 3335.182
$ws.Range("I113").Value = 2841
$ws.Range("J113").Value = 4200
$ws.Range("K113").Value = 2841
$ws.Range("L113").Value = 4200
$ws.Range("M113").Value = 413
$ws.Range("N113").Value = -10708
$ws.Range("H116").Value = 15759.706
$ws.Range("J116").Value = 14000
$ws.Range("L116").Value = 14000
$ws.Range("N116").Value = -20884
$ws.Range("H136").Value = 60500
$ws.Range("J136").Value = 60500
$ws.Range("L136").Value = 60500
$ws.Range("N136").Value = -70700
$ws.Range("H137").Value = 2135.5334
$ws.Range("I137").Value = 2080
$ws.Range("J137").Value = 2199
$ws.Range("K137").Value = 6240
$ws.Range("L137").Value = 6597
$ws.Range("M137").Value = -3690
$ws.Range("N137").Value = -11697

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 8328.666999999999
$ws.Range("I6").Value = 8328.666999999999
$ws.Range("K6").Value = 8328.666999999999
$ws.Range("M6").Value = -8155.666999999999
$ws.Range("H32").Value = 8595.339
$ws.Range("I32").Value = 5380.875
$ws.Range("J32").Value = 19616.357
$ws.Range("K32").Value = 5380.875
$ws.Range("L32").Value = 19616.357
$ws.Range("M32").Value = -5093.875
$ws.Range("N32").Value = -20190.357
$ws.Range("H45").Value = 3573.2354
$ws.Range("I45").Value = 3269.2222
$ws.Range("J45").Value = 3915.25
$ws.Range("K45").Value = 3269.2222
$ws.Range("L45").Value = 3915.25
$ws.Range("M45").Value = -2892.2222
$ws.Range("N45").Value = -4669.25
$ws.Range("H74").Value = 6195.9414
$ws.Range("I74").Value = 3174.5
$ws.Range("J74").Value = 7125.615
$ws.Range("K74").Value = 3174.5
$ws.Range("L74").Value = 7125.615
$ws.Range("M74").Value = -2300.5
$ws.Range("N74").Value = -8873.615
$ws.Range("H77").Value = 6195.9414
$ws.Range("I77").Value = 3174.5
$ws.Range("J77").Value = 7125.615
$ws.Range("K77").Value = 15872.5
$ws.Range("L77").Value = 35628.075
$ws.Range("M77").Value = -11504.5
$ws.Range("N77").Value = -44364.075

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 57332.668
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 57332.668
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 57332.668
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -57962.668
$ws.Range("H79").Value = 57332.668
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 57332.668
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 57332.668
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -59516.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1600.2727
$ws.Range("I22").Value = 1490.3636
$ws.Range("K22").Value = 1490.3636
$ws.Range("M22").Value = -1140.3636
$ws.Range("H31").Value = 8987.674000000001
$ws.Range("I31").Value = 5111.6294
$ws.Range("K31").Value = 5111.6294
$ws.Range("M31").Value = -4816.6294
$ws.Range("H34").Value = 8987.674000000001
$ws.Range("I34").Value = 5111.6294
$ws.Range("K34").Value = 5111.6294
$ws.Range("M34").Value = -4909.6294
$ws.Range("H86").Value = 4784.294
$ws.Range("J86").Value = 4978.6665
$ws.Range("L86").Value = 4978.6665
$ws.Range("N86").Value = -7224.6665
$ws.Range("H89").Value = 4784.294
$ws.Range("J89").Value = 4978.6665
$ws.Range("L89").Value = 24893.3325
$ws.Range("N89").Value = -36125.3325
$ws.Range("H108").Value = 107149.6
$ws.Range("J108").Value = 107149.6
$ws.Range("L108").Value = 107149.6
$ws.Range("N108").Value = -114829.6
$ws.Range("H132").Value = 4605
$ws.Range("I132").Value = 2683.8845
$ws.Range("J132").Value = 9599.9
$ws.Range("K132").Value = 8051.6535
$ws.Range("L132").Value = 28799.7
$ws.Range("M132").Value = -5521.6535
$ws.Range("N132").Value = -33859.7
$ws.Range("H134").Value = 3441.1765
$ws.Range("I134").Value = 3400.125
$ws.Range("K134").Value = 10200.375
$ws.Range("M134").Value = -7665.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 515.9583
$ws.Range("J92").Value = 446.17648
$ws.Range("L92").Value = 1338.52944
$ws.Range("N92").Value = -3834.52944
$ws.Range("H97").Value = 3335.524
$ws.Range("I97").Value = 5253.091
$ws.Range("K97").Value = 15759.273
$ws.Range("M97").Value = -15263.273
$ws.Range("H105").Value = 20000
$ws.Range("J105").Value = 20000
$ws.Range("L105").Value = 60000
$ws.Range("N105").Value = -65242
$ws.Range("H139").Value = 2474.5557
$ws.Range("I139").Value = 2474.5557
$ws.Range("K139").Value = 7423.6671
$ws.Range("M139").Value = -2283.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5424.778
$ws.Range("I113").Value = 5442
$ws.Range("K113").Value = 5442
$ws.Range("M113").Value = -3272
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H138").Value = 99999.10000000001
$ws.Range("J138").Value = 99999.10000000001
$ws.Range("L138").Value = 99999.10000000001
$ws.Range("N138").Value = -110279.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3068.5625
$ws.Range("I46").Value = 1687
$ws.Range("J46").Value = 3897.5
$ws.Range("K46").Value = 1687
$ws.Range("L46").Value = 3897.5
$ws.Range("M46").Value = -1499
$ws.Range("N46").Value = -4273.5
$ws.Range("H61").Value = 4795.3076
$ws.Range("I61").Value = 1472.25
$ws.Range("J61").Value = 10112.2
$ws.Range("K61").Value = 1472.25
$ws.Range("L61").Value = 10112.2
$ws.Range("M61").Value = -1270.25
$ws.Range("N61").Value = -10516.2
$ws.Range("H113").Value = 4795.3076
$ws.Range("I113").Value = 1472.25
$ws.Range("J113").Value = 10112.2
$ws.Range("K113").Value = 1472.25
$ws.Range("L113").Value = 10112.2
$ws.Range("M113").Value = 697.75
$ws.Range("N113").Value = -14452.2
$ws.Range("H132").Value = 5875.8887
$ws.Range("I132").Value = 4422.875
$ws.Range("J132").Value = 17500
$ws.Range("K132").Value = 13268.625
$ws.Range("L132").Value = 52500
$ws.Range("M132").Value = -10738.625
$ws.Range("N132").Value = -57560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1209.0444
$ws.Range("I107").Value = 1291.3214
$ws.Range("K107").Value = 3873.9642
$ws.Range("M107").Value = -1953.9642
$ws.Range("H113").Value = 462.73334
$ws.Range("I113").Value = 452.07693
$ws.Range("K113").Value = 1356.23079
$ws.Range("M113").Value = 813.7692099999999
$ws.Range("H136").Value = 6132.679
$ws.Range("I136").Value = 3269.5405
$ws.Range("J136").Value = 12753.6875
$ws.Range("K136").Value = 9808.621500000001
$ws.Range("L136").Value = 38261.0625
$ws.Range("M136").Value = -7258.621500000001
$ws.Range("N136").Value = -43361.0625
